$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.160.57'
$ws.Range("E2").Value = '  -2.61%  '

$ws.Range("D3").Value = '1.869.22'
$ws.Range("E3").Value = '  -2.06%  '

$ws.Range("E4").Value = '  -0.28%  '

$ws.Range("D5").Value = '307.05'
$ws.Range("E5").Value = '  -2.07%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.21%  '

$ws.Range("D7").Value = '0.5155'
$ws.Range("E7").Value = '  +3.12%  '

$ws.Range("D8").Value = '0.3753'
$ws.Range("E8").Value = '  -1.72%  '

$ws.Range("D9").Value = '0.07175'
$ws.Range("E9").Value = '  -1.68%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '0.8872'
$ws.Range("E10").Value = '  -2.35%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = '20.70'
$ws.Range("E11").Value = '  -1.02%  '

$ws.Range("E12").Value = '  -1.28%  '

$ws.Range("D13").Value = '1.859.84'
$ws.Range("E13").Value = '  -3.20%  '

$ws.Range("E14").Value = '  -2.74%  '

$ws.Range("D15").Value = '89.43'
$ws.Range("E15").Value = '  -2.51%  '

$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.24%  '

$ws.Range("E17").Value = '  -1.86%  '

$ws.Range("D18").Value = '14.18'
$ws.Range("E18").Value = '  -2.84%  '

$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.30%  '

$ws.Range("D20").Value = '27.192.86'
$ws.Range("E20").Value = '  -2.61%  '

$ws.Range("D21").Value = '5.059'
$ws.Range("E21").Value = '  -2.41%  '

$ws.Range("D22").Value = '2.079.39'
$ws.Range("E22").Value = '  -3.86%  '

$ws.Range("D23").Value = '10.62'
$ws.Range("E23").Value = '  -1.96%  '

$ws.Range("D24").Value = '6.483'
$ws.Range("E24").Value = '  -1.79%  '

$ws.Range("D25").Value = '150.93'
$ws.Range("E25").Value = '  -2.20%  '

$ws.Range("E26").Value = '  -1.89%  '

$ws.Range("D27").Value = '18.02'
$ws.Range("E27").Value = '  -2.15%  '

$ws.Range("D28").Value = '2.125'
$ws.Range("E28").Value = '  -4.59%  '

$ws.Range("D29").Value = '112.74'
$ws.Range("E29").Value = '  -2.35%  '

$ws.Range("D30").Value = '4.762'
$ws.Range("E30").Value = '  -3.15%  '

$ws.Range("E31").Value = '  +0.64%  '

$ws.Range("D32").Value = '0.09002'
$ws.Range("E32").Value = '  +0.11%  '

$ws.Range("D33").Value = '0.05150'
$ws.Range("E33").Value = '  -2.09%  '

$ws.Range("D34").Value = '3.103'
$ws.Range("E34").Value = '  -3.33%  '

$ws.Range("D35").Value = '0.7528'
$ws.Range("E35").Value = '  -1.69%  '

$ws.Range("D36").Value = '1.171'
$ws.Range("E36").Value = '  -5.25%  '

$ws.Range("E37").Value = '  -1.23%  '

$ws.Range("D38").Value = '2.528'
$ws.Range("E38").Value = '  -0.85%  '

$ws.Range("D39").Value = '3.028'
$ws.Range("E39").Value = '  +0.14%  '

$ws.Range("D40").Value = '1.081'
$ws.Range("E40").Value = '  -1.11%  '

$ws.Range("D41").Value = '0.5351'
$ws.Range("E41").Value = '  -4.12%  '

$ws.Range("D42").Value = '6.646'

$ws.Range("D43").Value = '114.75'
$ws.Range("E43").Value = '  +2.91%  '

$ws.Range("E44").Value = '  -0.37%  '

$ws.Range("E45").Value = '  -1.93%  '

$ws.Range("D46").Value = '0.4669'
$ws.Range("E46").Value = '  -3.36%  '

$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  -0.32%  '

$ws.Range("E48").Value = '  -4.46%  '

$ws.Range("D49").Value = '1.573'
$ws.Range("E49").Value = '  -3.59%  '

$ws.Range("D50").Value = '65.01'
$ws.Range("E50").Value = '  -3.89%  '

$ws.Range("D51").Value = '36.47'
$ws.Range("E51").Value = '  -1.58%  '

